# Commit: "added REST_API + Database"
#
# The presentation's slide deck is trimmed down: the three interim
# "diagram" slides (ER-Diagramm, Ablaufdiagramm, Klassendiagramm) that sat
# between "Aktueller Stand" and "Nächste Schritte" are removed, leaving
# "Nächste Schritte" (which already talks about the REST API / database
# work) directly after "Aktueller Stand".
#
# Before: LogIt, Beschreibung, Aktueller Stand, ER-Diagramm,
#         Ablaufdiagramm, Klassendiagramm, Nächste Schritte
# After:  LogIt, Beschreibung, Aktueller Stand, Nächste Schritte

$p = $ppt.ActivePresentation

# The three diagram slides all sit at slide index 4 once the previous one
# is removed (Aktueller Stand = 3, then the next unwanted slide always
# slides into position 4).
$p.Slides.Item(4).Delete()
$p.Slides.Item(4).Delete()
$p.Slides.Item(4).Delete()
